$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cells that Excel will not mis-parse as numbers/dates.
$ws.Range("A1").Value = "Henry"
$ws.Range("B1").Value = "Hunter"
$ws.Range("C1").Value = "TE"

# Numeric-looking values that must stay as literal text -- force text
# interpretation via NumberFormat, then drop the format again so no
# stray cell style is left behind.
$ws.Range("D1:F1").NumberFormat = "@"
$ws.Range("D1").Value = "2019-01-13"
$ws.Range("E1").Value = "18"
$ws.Range("F1").Value = "24.037"
$ws.Range("D1:F1").ClearFormats()

$ws.Range("G1").Value = "LAC"
$ws.Range("H1").Value = "@"
$ws.Range("I1").Value = "NWE"
$ws.Range("J1").Value = "L 28-41"
$ws.Range("K1").Value = "*"

# Trailing numeric cell.
$ws.Range("L1").Value = 0
